$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, pushing existing rows 19-80 down to 20-81.
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new weekly price entry.
$ws.Cells.Item(19, 1).Value2 = 10
$ws.Cells.Item(19, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(19, 3).Value2 = "La Araucanía"
$ws.Cells.Item(19, 4).Value2 = 45071
$ws.Cells.Item(19, 5).Value2 = 9
$ws.Cells.Item(19, 6).Value2 = 100112042
$ws.Cells.Item(19, 7).Value2 = "Locoto"
$ws.Cells.Item(19, 8).Value2 = "Sin especificar"
$ws.Cells.Item(19, 9).Value2 = "Primera"
$ws.Cells.Item(19, 10).Value2 = 150
$ws.Cells.Item(19, 11).Value2 = 4400
$ws.Cells.Item(19, 12).Value2 = 4400
$ws.Cells.Item(19, 13).Value2 = 4400
$ws.Cells.Item(19, 14).Value2 = "`$/kilo"
$ws.Cells.Item(19, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(19, 16).Value2 = 4400
$ws.Cells.Item(19, 17).Value2 = 1
$ws.Cells.Item(19, 18).Value2 = "Hortaliza"
